# Update crypto price/volume data per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cell -> new text value. All D/E columns are stored as plain text
# (inline strings) in the sheet, so every value is written as text; a leading
# apostrophe is used for values that would otherwise be auto-parsed by Excel
# as a number (e.g. "403.36"), then ClearFormats() strips the resulting
# quote-prefix formatting so no stray number-format/style is left behind.
$updates = [ordered]@{
    "D2" = "65.943.27"
    "E2" = "  -1.98%  "
    "D3" = "3.768.45"
    "E3" = "  +1.44%  "
    "E4" = "  +0.11%  "
    "D5" = "'403.36"
    "E5" = "  -4.84%  "
    "D6" = "'130.52"
    "E6" = "  -0.77%  "
    "D7" = "3.758.51"
    "E7" = "  +1.37%  "
    "D8" = "'0.603"
    "E8" = "  -6.28%  "
    "E9" = "  +0.00%  "
    "D10" = "'0.720"
    "E10" = "  -6.34%  "
    "D11" = "'0.165"
    "E11" = "  -10.59%  "
    "D12" = "'0.0000355"
    "E12" = "  -10.32%  "
    "E13" = "  -5.91%  "
    "D14" = "4.369.49"
    "E14" = "  +1.63%  "
    "D15" = "'9.65"
    "E15" = "  -5.01%  "
    "D16" = "'14.58"
    "E16" = "  +12.70%  "
    "E17" = "  -1.55%  "
    "D18" = "3.774.52"
    "E18" = "  +0.85%  "
    "D19" = "'19.33"
    "E19" = "  -7.09%  "
    "D20" = "66.368.42"
    "E20" = "  -1.28%  "
    "E21" = "  -6.68%  "
    "D22" = "'410.86"
    "E22" = "  -8.97%  "
    "E23" = "  -9.18%  "
    "D24" = "'84.55"
    "E24" = "  -5.49%  "
    "E25" = "  -5.36%  "
    "B26" = "LEO"
    "C26" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D26" = "'5.69"
    "E26" = "  +14.33%  "
    "B27" = "EthereumClassic"
    "C27" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D27" = "'36.10"
    "E27" = "  -5.24%  "
    "D28" = "'3.08"
    "E28" = "  -7.39%  "
    "D29" = "'9.26"
    "E29" = "  -9.67%  "
    "D30" = "'12.28"
    "E30" = "  -3.42%  "
    "D31" = "'2.72"
    "E31" = "  -2.48%  "
    "D32" = "'0.117"
    "E32" = "  -4.37%  "
    "D33" = "'7.21"
    "E33" = "  -1.84%  "
    "E34" = "  -6.17%  "
    "D35" = "'38.74"
    "E35" = "  -8.14%  "
    "D36" = "'0.999"
    "E36" = "  -0.02%  "
    "D37" = "'55.08"
    "E37" = "  -2.41%  "
    "D38" = "0.0₃0730"
    "E38" = "  -5.53%  "
    "D39" = "'0.0455"
    "E39" = "  -7.65%  "
    "D40" = "'2.86"
    "E40" = "  -9.95%  "
    "D41" = "'0.999"
    "E41" = "  +0.14%  "
    "D42" = "'0.134"
    "E42" = "  -8.65%  "
    "D43" = "'27.09"
    "E43" = "  -4.06%  "
    "B44" = "ApeXProtocol"
    "C44" = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
    "D44" = "'3.13"
    "E44" = "  +17.97%  "
    "B45" = "Monero"
    "C45" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D45" = "'144.92"
    "E45" = "  -1.09%  "
    "D46" = "'3.23"
    "E46" = "  -6.12%  "
    "E47" = "  -3.64%  "
    "B48" = "WEMIXToken"
    "C48" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D48" = "'2.56"
    "E48" = "  -4.69%  "
    "B49" = "NEARProtocol"
    "C49" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D49" = "'4.21"
    "E49" = "  -5.19%  "
    "B50" = "Stacks"
    "C50" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D50" = "'2.78"
    "E50" = "  -5.80%  "
    "D51" = "'0.289"
    "E51" = "  -6.65%  "
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    $range.Value = $updates[$cellRef]
    $range.ClearFormats()
}
